# Update countries & provincias Spain
#
# This applies the "2020-08-12 01:47 -> 03:04" refresh of the COVID-19
# country table: refreshed figures for several rows, plus three pairs of
# adjacent countries that swapped places in the sorted order (their figures
# travel with the country name).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- timestamp caption (A1) ----------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 12 de Agosto de 2020 a las 03:04"

# --- plain figure refreshes (country stays on the same row) --------------

# Row 4: Estados Unidos
$ws.Range("B4").Value = 5304378
$ws.Range("C4").Value = 52940
$ws.Range("D4").Value = 2755089
$ws.Range("E4").Value = 2381700
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1397
$ws.Range("H4").Value = 167589

# Row 15: Reino Unido
$ws.Range("G15").Value = 102
$ws.Range("H15").Value = 46628

# Row 19: Argentina
$ws.Range("B19").Value = 260911
$ws.Range("C19").Value = 7043
$ws.Range("E19").Value = 74518
$ws.Range("G19").Value = 240
$ws.Range("H19").Value = 5004

# Row 39: Panama
$ws.Range("B39").Value = 76464
$ws.Range("C39").Value = 1070
$ws.Range("D39").Value = 50665
$ws.Range("E39").Value = 24119
$ws.Range("G39").Value = 16
$ws.Range("H39").Value = 1680

# Row 159: Vietnam
$ws.Range("B159").Value = 866
$ws.Range("C159").Value = 19
$ws.Range("D159").Value = 451
$ws.Range("E159").Value = 399

# Row 181: Papua Nueva Guinea
$ws.Range("D181").Value = 71
$ws.Range("E181").Value = 140

# Row 195: San Martin (Parte Francesa)
$ws.Range("E195").Value = 36
$ws.Range("G195").Value = 1
$ws.Range("H195").Value = 4

# --- reorderings: country name + its whole data row move together --------

# Rows 65-67 (between Moldavia@64 and Costa Rica@68):
# before: Kenia, Irlanda, Venezuela
# after : Venezuela (new figures), Kenia (old row65 figures), Irlanda (old row66 figures)
$ws.Range("A65").Value = "Venezuela"
$ws.Range("B65").Value = 27938
$ws.Range("C65").Value = 1138
$ws.Range("D65").Value = 19706
$ws.Range("E65").Value = 7994
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 9
$ws.Range("H65").Value = 238

$ws.Range("A66").Value = "Kenia"
$ws.Range("B66").Value = 27425
$ws.Range("C66").Value = 497
$ws.Range("D66").Value = 13867
$ws.Range("E66").Value = 13120
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 15
$ws.Range("H66").Value = 438

$ws.Range("A67").Value = "Irlanda"
$ws.Range("B67").Value = 26801
$ws.Range("C67").Value = 33
$ws.Range("D67").Value = 23364
$ws.Range("E67").Value = 1664
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 1773

# Rows 101-102 (between Mauritania@100 and Croacia@103):
# before: Grecia, Libia
# after : Libia (new figures), Grecia (old row101 figures)
$ws.Range("A101").Value = "Libia"
$ws.Range("B101").Value = 6302
$ws.Range("C101").Value = 373
$ws.Range("D101").Value = 740
$ws.Range("E101").Value = 5430
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 7
$ws.Range("H101").Value = 132

$ws.Range("A102").Value = "Grecia"
$ws.Range("B102").Value = 5942
$ws.Range("C102").Value = 193
$ws.Range("D102").Value = 3804
$ws.Range("E102").Value = 1924
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 214

# Rows 213-214 (between Bonaire...@212 and Santa Sede@215):
# before: Montserrat, Islas Malvinas
# after : Islas Malvinas (old row214 figures), Montserrat (old row213 figures)
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 13
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1
